$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.363.50"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.15%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.654.93"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -0.07%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "519.13"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.18%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.39"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.20%  "
$ws.Range("E7").Value = "  +0.25%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.569"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.59%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.663.78"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.10%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.92"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +9.05%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.102"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -2.67%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.334"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -1.04%  "
$ws.Range("E13").Value = "  +1.82%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.127.98"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.30%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "59.422.15"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.26%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.00"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.09%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.673.12"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.42%  "
$ws.Range("E18").Value = "  -1.77%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "339.17"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -4.00%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.42"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.83%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.31"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -1.11%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.29"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.88%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.993"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.54%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.28"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +2.16%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.167"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.51%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.411"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -1.70%  "
$ws.Range("E27").Value = "  +0.59%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0801"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -1.37%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.14"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.06%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.68"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +5.14%  "
$ws.Range("E31").Value = "  +0.09%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.58"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.47%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.73"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -1.41%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "149.49"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.17%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.15"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +1.75%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.20"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.68%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.901"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -5.75%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.879"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +1.46%  "
$ws.Range("B39").Value = "OKB"
$ws.Range("C39").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.92"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.68%  "
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.49"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +3.39%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.57"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -2.97%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.629"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +3.68%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.999"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.17%  "
$ws.Range("B44").Value = "Bittensor"
$ws.Range("C44").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "275.33"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -1.99%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "19.78"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.41%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0975"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -1.62%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0534"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.79%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.056.08"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -3.23%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "10.51"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +2.11%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.78"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -1.29%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0229"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -1.44%  "
